$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106; this shifts existing rows 106:242 down to 107:243
$ws.Rows(106).Insert()

# Populate the newly inserted row 106 with the new data record
$ws.Range("A106").Value = 3
$ws.Range("B106").Value = "Femacal de La Calera"
$ws.Range("C106").Value = "Coquimbo"
$ws.Range("D106").Value = 44546
$ws.Range("E106").Value = 5
$ws.Range("F106").Value = 100112012
$ws.Range("G106").Value = "Espinaca"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 188
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = 3000
$ws.Range("N106").Value = "$/docena de atados (3 kilos)"
$ws.Range("O106").Value = "Provincia de Quillota"
$ws.Range("P106").Value = 1000
$ws.Range("Q106").Value = 3
$ws.Range("R106").Value = "Hortaliza"
